$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; Col='D'; Value='28.385.25'; ForceText=$false},
    @{Row=2; Col='E'; Value='  -2.71%  '; ForceText=$false},
    @{Row=3; Col='D'; Value='1.952.48'; ForceText=$false},
    @{Row=3; Col='E'; Value='  -0.66%  '; ForceText=$false},
    @{Row=4; Col='D'; Value='1.006'; ForceText=$true},
    @{Row=4; Col='E'; Value='  -1.45%  '; ForceText=$false},
    @{Row=5; Col='D'; Value='319.96'; ForceText=$true},
    @{Row=5; Col='E'; Value='  -2.45%  '; ForceText=$false},
    @{Row=6; Col='E'; Value='  -1.22%  '; ForceText=$false},
    @{Row=7; Col='D'; Value='0.4762'; ForceText=$true},
    @{Row=7; Col='E'; Value='  -4.45%  '; ForceText=$false},
    @{Row=8; Col='D'; Value='0.4027'; ForceText=$true},
    @{Row=8; Col='E'; Value='  -3.90%  '; ForceText=$false},
    @{Row=9; Col='D'; Value='53.52'; ForceText=$true},
    @{Row=9; Col='E'; Value='  +0.07%  '; ForceText=$false},
    @{Row=10; Col='D'; Value='0.08406'; ForceText=$true},
    @{Row=10; Col='E'; Value='  -5.01%  '; ForceText=$false},
    @{Row=11; Col='D'; Value='1.054'; ForceText=$true},
    @{Row=11; Col='E'; Value='  -3.86%  '; ForceText=$false},
    @{Row=12; Col='E'; Value='  -3.03%  '; ForceText=$false},
    @{Row=13; Col='D'; Value='1.953.28'; ForceText=$false},
    @{Row=13; Col='E'; Value='  -8.13%  '; ForceText=$false},
    @{Row=14; Col='D'; Value='7.550'; ForceText=$true},
    @{Row=14; Col='E'; Value='  -4.00%  '; ForceText=$false},
    @{Row=15; Col='D'; Value='6.134'; ForceText=$true},
    @{Row=15; Col='E'; Value='  -4.10%  '; ForceText=$false},
    @{Row=16; Col='D'; Value='1.008'; ForceText=$true},
    @{Row=16; Col='E'; Value='  -1.16%  '; ForceText=$false},
    @{Row=17; Col='D'; Value='90.50'; ForceText=$true},
    @{Row=17; Col='E'; Value='  -0.70%  '; ForceText=$false},
    @{Row=18; Col='D'; Value='0.00001066'; ForceText=$true},
    @{Row=18; Col='E'; Value='  -2.96%  '; ForceText=$false},
    @{Row=19; Col='D'; Value='0.06577'; ForceText=$true},
    @{Row=19; Col='E'; Value='  -2.16%  '; ForceText=$false},
    @{Row=20; Col='D'; Value='18.46'; ForceText=$true},
    @{Row=20; Col='E'; Value='  -3.72%  '; ForceText=$false},
    @{Row=21; Col='D'; Value='1.006'; ForceText=$true},
    @{Row=21; Col='E'; Value='  -1.52%  '; ForceText=$false},
    @{Row=22; Col='D'; Value='5.816'; ForceText=$true},
    @{Row=22; Col='E'; Value='  -1.94%  '; ForceText=$false},
    @{Row=23; Col='D'; Value='28.398.90'; ForceText=$false},
    @{Row=23; Col='E'; Value='  -2.98%  '; ForceText=$false},
    @{Row=24; Col='D'; Value='11.44'; ForceText=$true},
    @{Row=24; Col='E'; Value='  -3.60%  '; ForceText=$false},
    @{Row=25; Col='D'; Value='2.287'; ForceText=$true},
    @{Row=25; Col='E'; Value='  -1.11%  '; ForceText=$false},
    @{Row=26; Col='D'; Value='2.194.13'; ForceText=$false},
    @{Row=26; Col='E'; Value='  -10.57%  '; ForceText=$false},
    @{Row=27; Col='D'; Value='154.61'; ForceText=$true},
    @{Row=27; Col='E'; Value='  -0.83%  '; ForceText=$false},
    @{Row=28; Col='D'; Value='20.15'; ForceText=$true},
    @{Row=28; Col='E'; Value='  -2.10%  '; ForceText=$false},
    @{Row=29; Col='E'; Value='  -4.56%  '; ForceText=$false},
    @{Row=30; Col='D'; Value='2.148'; ForceText=$true},
    @{Row=30; Col='E'; Value='  -5.83%  '; ForceText=$false},
    @{Row=31; Col='D'; Value='123.37'; ForceText=$true},
    @{Row=31; Col='E'; Value='  -2.37%  '; ForceText=$false},
    @{Row=32; Col='D'; Value='0.9766'; ForceText=$true},
    @{Row=32; Col='E'; Value='  -6.39%  '; ForceText=$false},
    @{Row=33; Col='E'; Value='  -2.63%  '; ForceText=$false},
    @{Row=34; Col='D'; Value='1.448'; ForceText=$true},
    @{Row=34; Col='E'; Value='  -2.97%  '; ForceText=$false},
    @{Row=35; Col='E'; Value='  -2.01%  '; ForceText=$false},
    @{Row=36; Col='E'; Value='  -3.34%  '; ForceText=$false},
    @{Row=37; Col='D'; Value='8.954'; ForceText=$true},
    @{Row=37; Col='E'; Value='  -2.17%  '; ForceText=$false},
    @{Row=38; Col='D'; Value='0.02323'; ForceText=$true},
    @{Row=38; Col='E'; Value='  -4.34%  '; ForceText=$false},
    @{Row=39; Col='D'; Value='0.06209'; ForceText=$true},
    @{Row=39; Col='E'; Value='  -1.47%  '; ForceText=$false},
    @{Row=40; Col='D'; Value='1.244'; ForceText=$true},
    @{Row=40; Col='E'; Value='  -3.27%  '; ForceText=$false},
    @{Row=41; Col='D'; Value='0.6195'; ForceText=$true},
    @{Row=41; Col='E'; Value='  -4.03%  '; ForceText=$false},
    @{Row=42; Col='E'; Value='  -3.69%  '; ForceText=$false},
    @{Row=43; Col='E'; Value='  -1.18%  '; ForceText=$false},
    @{Row=44; Col='D'; Value='0.1916'; ForceText=$true},
    @{Row=44; Col='E'; Value='  -4.88%  '; ForceText=$false},
    @{Row=45; Col='D'; Value='1.341'; ForceText=$true},
    @{Row=45; Col='E'; Value='  +5.32%  '; ForceText=$false},
    @{Row=46; Col='B'; Value='Decentraland'; ForceText=$false},
    @{Row=46; Col='C'; Value='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; ForceText=$false},
    @{Row=46; Col='D'; Value='0.5943'; ForceText=$true},
    @{Row=46; Col='E'; Value='  -4.59%  '; ForceText=$false},
    @{Row=47; Col='B'; Value='EnergySwap'; ForceText=$false},
    @{Row=47; Col='C'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText=$false},
    @{Row=47; Col='D'; Value='12.92'; ForceText=$true},
    @{Row=47; Col='E'; Value='  -4.22%  '; ForceText=$false},
    @{Row=48; Col='D'; Value='2.051'; ForceText=$true},
    @{Row=48; Col='E'; Value='  -5.91%  '; ForceText=$false},
    @{Row=49; Col='D'; Value='3.386'; ForceText=$true},
    @{Row=49; Col='E'; Value='  -2.74%  '; ForceText=$false},
    @{Row=50; Col='D'; Value='0.00000000320'; ForceText=$true},
    @{Row=50; Col='E'; Value='  -4.03%  '; ForceText=$false},
    @{Row=51; Col='D'; Value='0.06799'; ForceText=$true},
    @{Row=51; Col='E'; Value='  -1.35%  '; ForceText=$false}

)

foreach ($chg in $changes) {
    $cellRef = "{0}{1}" -f $chg.Col, $chg.Row
    $rng = $ws.Range($cellRef)
    if ($chg.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}
